$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the lag-window rows (7-10) and refresh their AIC_sum values
$ws.Range("B7").Value = "1-17 lag"
$ws.Range("C7").Value = 441323.4472703643

$ws.Range("B8").Value = "1-19 lag"
$ws.Range("C8").Value = 440942.54606649885

$ws.Range("B9").Value = "1-23 lag"
$ws.Range("C9").Value = 440264.828535935

$ws.Range("B10").Value = "1-25 lag"
$ws.Range("C10").Value = 439907.7146539751

# Add the new "Population as offset" row (was row 10, now row 11)
# Force A11 to be stored as text "10" (matching the A column's text style)
# rather than the number 10: write a text formula, then paste-special as
# values so the literal text sticks without touching cell styles.
$ws.Range("A11").Formula = "=""10"""
$ws.Range("A11").Copy()
$ws.Range("A11").PasteSpecial(-4163)
$ws.Range("B11").Value = "Population as offset"
$ws.Range("C11").Value = 440364.74198841676
